$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.216.20'
$ws.Range("E2").Value = '  +2.06%  '
$ws.Range("D3").Value = '2.701.87'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'616.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.72%  '
$ws.Range("D6").Value = "'159.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.04%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("E9").Value = '  +7.22%  '
$ws.Range("D10").Value = "'6.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.44%  '
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").Value = "'0.0000211"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.36%  '
$ws.Range("D14").Value = "'30.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.77%  '
$ws.Range("D15").Value = '3.185.59'
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").Value = '66.050.67'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("D17").Value = '2.696.78'
$ws.Range("E17").Value = '  +3.48%  '
$ws.Range("D18").Value = "'12.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.86%  '
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").Value = "'7.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.23%  '
$ws.Range("D21").Value = "'360.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").Value = "'71.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.09%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = "'0.0000115"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +22.09%  '
$ws.Range("E25").Value = '  +6.46%  '
$ws.Range("D26").Value = "'1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").Value = "'1.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.08%  '
$ws.Range("E28").Value = '  +5.39%  '
$ws.Range("D29").Value = "'8.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.54%  '
$ws.Range("D30").Value = "'2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.83%  '
$ws.Range("D31").Value = "'541.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.18%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = "'1.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("E34").Value = '  +7.13%  '
$ws.Range("D35").Value = "'5.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").Value = "'20.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.77%  '
$ws.Range("D38").Value = "'164.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.24%  '
$ws.Range("D39").Value = "'2.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").Value = "'170.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.77%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = "'42.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.66%  '
$ws.Range("D44").Value = "'4.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.15%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = "'0.0628"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = "'2.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.19%  '
$ws.Range("D47").Value = "'23.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.04%  '
$ws.Range("E48").Value = '  +3.33%  '
$ws.Range("D49").Value = "'0.661"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("D50").Value = "'21.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.95%  '
$ws.Range("D51").Value = "'0.0997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.06%  '
